$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "/r/EarthPorn/comments/gz78ba/open_letter_to_steve_huffman_and_the_board_of/",
    "/r/EarthPorn/comments/hd5a2q/seascape_sunset_at_the_cliff_in_sydney_australia/",
    "/r/goodnews/comments/gwtp43/whats_new_content_creators_june_2020/",
    "/r/goodnews/comments/hc0d3l/good_news_its_friday_whats_your_feelgood_story/",
    "/r/learnpython/comments/h9527c/ask_anything_monday_weekly_thread/",
    "/r/learnpython/comments/hd6uou/book_recommendations_taking_udemy_course_and/",
    "/r/pics/comments/hd8660/michael_and_barbara_shetterly_listen_to_their_son/",
    "/r/Python/comments/gdfaip/rpython_job_board_for_may_june_july/",
    "/r/Python/comments/ha4qfy/whats_everyone_working_on_this_week/",
    "/r/EarthPorn/comments/hdd0k7/anyone_else_need_some_color_today_little_ohio/",
    "/r/learnpython/comments/hd9oaf/what_are_some_really_introbeginner_basic_project/",
    "/r/pics/comments/hda7w3/dubai_solar_eclipse_by_josh_cripps/"
)

$startRow = 12
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
